$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set C2 (name) before A2 (email) so the shared-string table is rebuilt
# in the same order as the target workbook: tester, Gopi, gopi@testleaf.com
$ws.Range("C2").Value = "Gopi"
$ws.Range("A2").Value = "gopi@testleaf.com"

# B2 keeps its existing value "tester" (index only shifts because of the
# shared string table changes above)

# Move the active selection to B2, matching the saved sheet view state
$ws.Range("B2").Select()
